$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.158.47"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").Value = "2.321.04"
$ws.Range("E3").Value = "  +0.62%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'303.55"
$ws.Range("E5").Value = "  +0.71%  "
$ws.Range("D6").Value = "'99.95"
$ws.Range("E6").Value = "  +0.50%  "
$ws.Range("E7").Value = "  +0.39%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "'0.519"
$ws.Range("E9").Value = "  +2.82%  "
$ws.Range("D10").Value = "'36.28"
$ws.Range("E10").Value = "  +6.00%  "
$ws.Range("D11").Value = "'0.0793"
$ws.Range("E11").Value = "  -0.57%  "
$ws.Range("E12").Value = "  -0.98%  "
$ws.Range("D13").Value = "'17.75"
$ws.Range("E13").Value = "  -0.77%  "
$ws.Range("D14").Value = "'6.92"
$ws.Range("E14").Value = "  +1.74%  "
$ws.Range("D15").Value = "2.683.19"
$ws.Range("E15").Value = "  +0.74%  "
$ws.Range("D16").Value = "2.377.88"
$ws.Range("E16").Value = "  +3.88%  "
$ws.Range("D17").Value = "'0.797"
$ws.Range("E17").Value = "  -2.01%  "
$ws.Range("D18").Value = "43.075.96"
$ws.Range("E18").Value = "  +0.10%  "
$ws.Range("D19").Value = "'13.11"
$ws.Range("E19").Value = "  +4.31%  "
$ws.Range("D20").Value = "'6.21"
$ws.Range("E20").Value = "  +1.98%  "
$ws.Range("E21").Value = "  +0.58%  "
$ws.Range("D22").Value = "'68.30"
$ws.Range("E22").Value = "  +0.83%  "
$ws.Range("D23").Value = "'240.46"
$ws.Range("E23").Value = "  +1.52%  "
$ws.Range("E24").Value = "  -2.00%  "
$ws.Range("E25").Value = "  +0.15%  "
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("D27").Value = "'25.58"
$ws.Range("E27").Value = "  +3.34%  "
$ws.Range("D28").Value = "'168.86"
$ws.Range("E28").Value = "  +0.57%  "
$ws.Range("D29").Value = "'34.16"
$ws.Range("E29").Value = "  +0.30%  "
$ws.Range("D30").Value = "'9.20"
$ws.Range("E30").Value = "  +0.65%  "
$ws.Range("E31").Value = "  -1.90%  "
$ws.Range("D32").Value = "'4.94"
$ws.Range("E32").Value = "  +8.71%  "
$ws.Range("D33").Value = "'5.17"
$ws.Range("E33").Value = "  +2.40%  "
$ws.Range("E34").Value = "  -0.10%  "
$ws.Range("D35").Value = "'17.94"
$ws.Range("E35").Value = "  +5.84%  "
$ws.Range("E36").Value = "  -1.60%  "
$ws.Range("E37").Value = "  +1.39%  "
$ws.Range("E38").Value = "  +2.05%  "
$ws.Range("E39").Value = "  +0.19%  "
$ws.Range("D40").Value = "'2.78"
$ws.Range("E40").Value = "  -0.85%  "
$ws.Range("E41").Value = "  +0.34%  "
$ws.Range("D42").Value = "1.995.13"
$ws.Range("E42").Value = "  +0.15%  "
$ws.Range("E43").Value = "  +1.60%  "
$ws.Range("D44").Value = "'2.25"
$ws.Range("E44").Value = "  -4.35%  "
$ws.Range("E45").Value = "  +1.37%  "
$ws.Range("D46").Value = "'17.73"
$ws.Range("E46").Value = "  +0.24%  "
$ws.Range("E47").Value = "  +0.51%  "
$ws.Range("D48").Value = "'55.11"
$ws.Range("E48").Value = "  -2.53%  "
$ws.Range("D49").Value = "'76.17"
$ws.Range("E49").Value = "  +8.46%  "
$ws.Range("D50").Value = "2.549.44"
$ws.Range("E50").Value = "  +0.69%  "
$ws.Range("D51").Value = "'1.56"
$ws.Range("E51").Value = "  +1.42%  "
